# s2cDNASample - further cleaning to metadata
# 1. Shared string "E7760" -> "E7420" (affects every G2:G39 cell that shares it)
# 2. G2:G39 gets a new font (Arial, 11pt, black) with General number format
# 3. H2:H39 gets an explicit =FALSE() formula (value stays FALSE/0), style unchanged
# 4. The sheet's active selection moves from H2:H39 to G2:G39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 39

# --- 1. Update the sample-number text in column G ---------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7420"
}

# --- 2. Re-style column G (Arial 11, black, General format) -----------------
# Build the target formatting once on G2, then fan it out with a format-only
# paste so every cell lands on the very same style entry instead of each one
# churning through its own chain of intermediate styles.
$template = $ws.Cells.Item($firstRow, 7)
$template.NumberFormat = "General"
$template.Font.Size = 11
$template.Font.Name = "Arial"
$template.Font.Color = 0

$template.Copy()
$ws.Range("G3:G39").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Column H keeps its style, just gains an explicit FALSE() formula ----
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# --- 4. Update the sheet selection to G2:G39 ---------------------------------
$ws.Range("G2:G39").Select()
